$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($range, $value)
    # Force the cell to stay text even when the literal looks numeric
    # (matches "Format Cells -> Text" then typing the value in real Excel).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

function Set-PlainValue {
    param($range, $value)
    $range.Value = $value
}

# --- Fill in the "Punto obtenido ahora julia" (E) and "Valor Objetivo ... Julia" (F)
# --- columns, plus fill out the previously-empty C/D columns for the last three problems.

# MitsosBarton2006Ex312 (rows 3-4)
Set-PlainValue $ws.Range("E3") "(3.45,1.85)"
Set-TextValue  $ws.Range("F3") "37.05"
Set-PlainValue $ws.Range("E4") "(3.45,1.85)"
Set-TextValue  $ws.Range("F4") "37.05"

# MitsosBarton2006Ex313 (rows 5-6)
Set-PlainValue $ws.Range("E5") "(2.3,4.45)"
Set-TextValue  $ws.Range("F5") " -2.16"
Set-PlainValue $ws.Range("E6") "(2.3,4.45)"
Set-TextValue  $ws.Range("F6") " -2.16"

# MitsosBarton2006Ex314 (rows 7-8)
Set-PlainValue $ws.Range("C7") "(2.1,3.3)"
Set-TextValue  $ws.Range("D7") "14.31"
Set-PlainValue $ws.Range("E7") "(2.1,-1.45)"
Set-TextValue  $ws.Range("F7") "5.52"
Set-PlainValue $ws.Range("C8") "(2.1,3.3)"
Set-TextValue  $ws.Range("D8") "14.31"
Set-PlainValue $ws.Range("E8") "(2.1,-1.45)"
Set-TextValue  $ws.Range("F8") "5.52"

# MitsosBarton2006Ex323 (rows 9-10)
Set-PlainValue $ws.Range("C9") "(2.89,0.3)"
Set-TextValue  $ws.Range("D9") "8.35"
Set-PlainValue $ws.Range("E9") "(2.89,0.3)"
Set-TextValue  $ws.Range("F9") "8.35"
Set-PlainValue $ws.Range("C10") "(2.89,0.3)"
Set-TextValue  $ws.Range("D10") "8.35"
Set-PlainValue $ws.Range("E10") "(2.89,0.3)"
Set-TextValue  $ws.Range("F10") "8.35"

# MorganPatrone2006a (rows 11-12)
Set-PlainValue $ws.Range("C11") "(4.5,2.8)"
Set-TextValue  $ws.Range("D11") " -7.3"
Set-PlainValue $ws.Range("E11") "(4.5,2.8)"
Set-TextValue  $ws.Range("F11") " -7.3"
Set-PlainValue $ws.Range("C12") "(4.5,2.8)"
Set-TextValue  $ws.Range("D12") " -7.3"
Set-PlainValue $ws.Range("E12") "(4.5,2.8)"
Set-TextValue  $ws.Range("F12") " -7.3"

# Column E is narrower now that it holds shorter values (~24.45 chars).
$ws.Columns.Item(5).ColumnWidth = 23.6363636363636

# Move the active selection.
$ws.Range("E4").Select() | Out-Null
